$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "43.075.23"
$ws.Range("E2").Value = "  +2.08%  "

# Row 3
$ws.Range("D3").Value = "2.309.34"
$ws.Range("E3").Value = "  +1.86%  "

# Row 4
$ws.Range("E4").Value = "  -0.01%  "

# Row 5
$ws.Range("D5").Value = "'302.63"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.99%  "

# Row 6
$ws.Range("D6").Value = "'102.00"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +6.18%  "

# Row 7
$ws.Range("E7").Value = "  +1.98%  "

# Row 8
$ws.Range("E8").Value = "  -0.06%  "

# Row 9
$ws.Range("E9").Value = "  +5.74%  "

# Row 10
$ws.Range("D10").Value = "'35.81"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +7.98%  "

# Row 11
$ws.Range("D11").Value = "'0.0797"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +1.17%  "

# Row 12
$ws.Range("E12").Value = "  +3.50%  "

# Row 13
$ws.Range("D13").Value = "'17.99"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +15.88%  "

# Row 14
$ws.Range("E14").Value = "  +4.10%  "

# Row 15
$ws.Range("D15").Value = "2.686.45"
$ws.Range("E15").Value = "  +2.49%  "

# Row 16
$ws.Range("D16").Value = "2.302.88"
$ws.Range("E16").Value = "  +0.61%  "

# Row 17
$ws.Range("D17").Value = "'0.815"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +4.22%  "

# Row 18
$ws.Range("D18").Value = "43.013.17"
$ws.Range("E18").Value = "  +2.11%  "

# Row 19
$ws.Range("D19").Value = "'12.72"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +9.00%  "

# Row 20
$ws.Range("D20").Value = "'6.19"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +3.46%  "

# Row 21
$ws.Range("D21").Value = "0.0₃0905"
$ws.Range("E21").Value = "  +1.76%  "

# Row 22
$ws.Range("D22").Value = "'67.96"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +2.39%  "

# Row 23
$ws.Range("D23").Value = "'237.68"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +1.21%  "

# Row 24
$ws.Range("D24").Value = "'2.22"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +12.86%  "

# Row 25
$ws.Range("D25").Value = "'2.47"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.53%  "

# Row 26
$ws.Range("E26").Value = "  -0.02%  "

# Row 27
$ws.Range("D27").Value = "'24.84"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +3.88%  "

# Row 28
$ws.Range("D28").Value = "'2.30"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +11.22%  "

# Row 29
$ws.Range("B29").Value = "InjectiveProtocol"
$ws.Range("C29").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D29").Value = "'34.67"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +3.08%  "

# Row 30
$ws.Range("B30").Value = "Monero"
$ws.Range("C30").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D30").Value = "'167.91"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -0.24%  "

# Row 31
$ws.Range("E31").Value = "  +0.94%  "

# Row 32
$ws.Range("D32").Value = "'1.00"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +0.05%  "

# Row 33
$ws.Range("E33").Value = "  +2.95%  "

# Row 34
$ws.Range("D34").Value = "'4.69"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +3.29%  "

# Row 35
$ws.Range("E35").Value = "  +3.80%  "

# Row 36
$ws.Range("D36").Value = "'17.10"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +3.31%  "

# Row 37
$ws.Range("D37").Value = "'0.0697"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +1.71%  "

# Row 38
$ws.Range("E38").Value = "  +3.47%  "

# Row 39
$ws.Range("D39").Value = "'2.84"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +2.01%  "

# Row 40
$ws.Range("E40").Value = "  +4.65%  "

# Row 41
$ws.Range("E41").Value = "  +1.46%  "

# Row 42
$ws.Range("D42").Value = "2.001.13"
$ws.Range("E42").Value = "  +1.56%  "

# Row 43
$ws.Range("E43").Value = "  -4.77%  "

# Row 44
$ws.Range("E44").Value = "  +4.40%  "

# Row 45
$ws.Range("D45").Value = "'10.30"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +7.74%  "

# Row 46
$ws.Range("D46").Value = "'17.75"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +1.71%  "

# Row 47
$ws.Range("D47").Value = "'2.88"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +3.81%  "

# Row 48
$ws.Range("D48").Value = "'56.20"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +7.65%  "

# Row 49
$ws.Range("D49").Value = "2.529.12"
$ws.Range("E49").Value = "  +1.41%  "

# Row 50
$ws.Range("E50").Value = "  +3.71%  "

# Row 51
$ws.Range("D51").Value = "'4.59"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +2.35%  "
